$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2; B="Bitcoin"; C="https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D="37.579.24"; E="  +1.51%  "},
    @{Row=3; B="Ethereum"; C="https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D="2.036.80"; E="  +2.57%  "},
    @{Row=4; B="TetherUSD"; C="https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D="1.00"; E="  -0.10%  "},
    @{Row=5; B="BNB"; C="https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D="257.92"; E="  +4.89%  "},
    @{Row=6; B="XRP"; C="https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D="0.624"; E="  -1.09%  "},
    @{Row=7; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="57.89"; E="  -5.77%  "},
    @{Row=8; B="USDC"; C="https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D="1.00"; E="  -0.01%  "},
    @{Row=9; B="Cardano"; C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D="0.386"; E="  +0.72%  "},
    @{Row=10; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="0.0797"; E="  -0.91%  "},
    @{Row=11; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.103"; E="  -1.55%  "},
    @{Row=12; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="14.86"; E="  -0.45%  "},
    @{Row=13; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="2.337.56"; E="  +2.60%  "},
    @{Row=14; B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="0.821"; E="  -3.29%  "},
    @{Row=15; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="21.44"; E="  -4.06%  "},
    @{Row=16; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="5.37"; E="  -2.18%  "},
    @{Row=17; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="2.043.42"; E="  +2.78%  "},
    @{Row=18; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="37.494.65"; E="  +1.48%  "},
    @{Row=19; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="70.12"; E="  -0.43%  "},
    @{Row=20; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.0₃0856"; E="  -1.05%  "},
    @{Row=21; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="5.22"; E="  +0.61%  "},
    @{Row=22; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="229.40"; E="  -0.48%  "},
    @{Row=23; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="2.66"; E="  +5.35%  "},
    @{Row=24; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="1.00"; E="  -0.03%  "},
    @{Row=25; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="2.34"; E="  -1.76%  "},
    @{Row=26; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="9.13"; E="  -2.08%  "},
    @{Row=27; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="163.71"; E="  -0.19%  "},
    @{Row=28; B="Kaspa"; C="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D="0.138"; E="  -6.02%  "},
    @{Row=29; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="20.07"; E="  +2.44%  "},
    @{Row=30; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="1.36"; E="  -0.68%  "},
    @{Row=31; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.120"; E="  -1.11%  "},
    @{Row=32; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.0666"; E="  +7.03%  "},
    @{Row=33; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="4.73"; E="  -2.86%  "},
    @{Row=34; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="4.51"; E="  -0.79%  "},
    @{Row=35; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="2.49"; E="  +8.84%  "},
    @{Row=36; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="3.46"; E="  +3.03%  "},
    @{Row=37; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="1.00"; E="  -0.10%  "},
    @{Row=38; B="WEMIXToken"; C="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D="1.81"; E="  +2.40%  "},
    @{Row=39; B="THORChain"; C="https://coinranking.com/coin/ybmU-kKU+thorchain-rune"; D="5.39"; E="  -3.06%  "},
    @{Row=40; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="3.03"; E="  +3.98%  "},
    @{Row=41; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.0968"; E="  -3.14%  "},
    @{Row=42; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.0217"; E="  +1.37%  "},
    @{Row=43; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="1.19"; E="  +0.34%  "},
    @{Row=44; B="Maker"; C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D="1.403.30"; E="  +1.75%  "},
    @{Row=45; B="InjectiveProtocol"; C="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D="16.17"; E="  -1.89%  "},
    @{Row=46; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="91.06"; E="  +0.82%  "},
    @{Row=47; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.05"; E="  +0.58%  "},
    @{Row=48; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="7.39"; E="  +1.08%  "},
    @{Row=49; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="2.87"; E="  +1.74%  "},
    @{Row=50; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="2.02"; E="  -0.93%  "},
    @{Row=51; B="RocketPoolETH"; C="https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"; D="2.227.64"; E="  +2.59%  "}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $dCell = $ws.Cells.Item($item.Row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $item.D
    $dCell.Style = "Normal"
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}

Write-Host "Applied crypto list update for" $data.Count "rows"
